$d = $word.ActiveDocument

$replacements = @(
    @("363÷7=", "981÷6="),
    @("537÷6=", "347÷7="),
    @("226÷8=", "332÷7="),
    @("368÷6=", "484÷6="),
    @("118÷2=", "235÷9="),
    @("829÷8=", "186÷3="),
    @("146÷2=", "401÷8="),
    @("868÷3=", "340÷7="),
    @("145÷2=", "116÷9="),
    @("138÷3=", "194÷2="),
    @("873÷8=", "165÷6="),
    @("116÷5=", "953÷3="),
    @("245÷4=", "353÷2="),
    @("584÷5=", "219÷9="),
    @("109÷8=", "124÷6="),
    @("147÷3=", "337÷5="),
    @("506÷9=", "637÷3="),
    @("704÷8=", "959÷2="),
    @("874÷3=", "688÷3="),
    @("745÷2=", "440÷9="),
    @("351÷6=", "290÷7="),
    @("607÷4=", "514÷7="),
    @("661÷4=", "937÷3="),
    @("853÷5=", "940÷4="),
    @("594÷3=", "628÷5=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
